# Update crypto market data cells per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.626.94'
$ws.Cells.Item(2, 5).Value = '  +1.10%  '

$ws.Cells.Item(3, 4).Value = '3.393.31'
$ws.Cells.Item(3, 5).Value = '  +0.24%  '

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

$ws.Cells.Item(5, 4).Value = "'577.02"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.98%  '

$ws.Cells.Item(6, 4).Value = "'143.19"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.74%  '

$ws.Cells.Item(7, 5).Value = '  +0.01%  '

$ws.Cells.Item(8, 5).Value = '  -0.34%  '

$ws.Cells.Item(9, 4).Value = "'7.61"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.06%  '

$ws.Cells.Item(10, 5).Value = '  -0.56%  '

$ws.Cells.Item(11, 5).Value = '  -0.93%  '

$ws.Cells.Item(12, 4).Value = '3.972.58'
$ws.Cells.Item(12, 5).Value = '  +0.20%  '

$ws.Cells.Item(13, 5).Value = '  -0.32%  '

$ws.Cells.Item(14, 4).Value = "'27.96"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.64%  '

$ws.Cells.Item(15, 4).Value = '3.394.96'
$ws.Cells.Item(15, 5).Value = '  +0.54%  '

$ws.Cells.Item(16, 5).Value = '  -1.08%  '

$ws.Cells.Item(17, 4).Value = '61.653.49'
$ws.Cells.Item(17, 5).Value = '  +0.94%  '

$ws.Cells.Item(18, 4).Value = "'6.12"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.69%  '

$ws.Cells.Item(19, 5).Value = '  +0.39%  '

$ws.Cells.Item(20, 4).Value = "'9.11"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.50%  '

$ws.Cells.Item(21, 4).Value = "'387.16"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.93%  '

$ws.Cells.Item(22, 5).Value = '  -0.60%  '

$ws.Cells.Item(23, 5).Value = '  -0.76%  '

$ws.Cells.Item(24, 4).Value = "'0.999"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.10%  '

$ws.Cells.Item(25, 5).Value = '  -2.64%  '

$ws.Cells.Item(26, 4).Value = "'0.182"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.01%  '

$ws.Cells.Item(27, 5).Value = '  -0.01%  '

$ws.Cells.Item(28, 5).Value = '  +1.27%  '

$ws.Cells.Item(29, 5).Value = '  -0.27%  '

$ws.Cells.Item(30, 5).Value = '  -0.59%  '

$ws.Cells.Item(31, 5).Value = '  +0.01%  '

$ws.Cells.Item(32, 5).Value = '  -0.32%  '

$ws.Cells.Item(33, 4).Value = "'23.34"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.30%  '

$ws.Cells.Item(34, 4).Value = "'6.93"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.48%  '

$ws.Cells.Item(35, 4).Value = "'169.25"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.43%  '

$ws.Cells.Item(36, 5).Value = '  +2.57%  '

$ws.Cells.Item(37, 4).Value = '3.424.80'
$ws.Cells.Item(37, 5).Value = '  +0.23%  '

$ws.Cells.Item(38, 5).Value = '  -0.03%  '

$ws.Cells.Item(39, 4).Value = "'27.48"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +3.21%  '

$ws.Cells.Item(40, 5).Value = '  -1.39%  '

$ws.Cells.Item(41, 5).Value = '  +0.37%  '

$ws.Cells.Item(42, 5).Value = '  +0.77%  '

$ws.Cells.Item(43, 4).Value = "'1.66"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.35%  '

$ws.Cells.Item(44, 5).Value = '  +2.69%  '

$ws.Cells.Item(45, 4).Value = '2.474.49'
$ws.Cells.Item(45, 5).Value = '  +0.66%  '

$ws.Cells.Item(46, 5).Value = '  -1.23%  '

$ws.Cells.Item(47, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(47, 4).Value = "'1.00"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.02%  '

$ws.Cells.Item(48, 2).Value = 'Cosmos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(48, 4).Value = "'6.61"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.64%  '

$ws.Cells.Item(49, 5).Value = '  -0.07%  '

$ws.Cells.Item(50, 4).Value = "'2.02"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -6.51%  '

$ws.Cells.Item(51, 5).Value = '  -1.43%  '
